# Upload fill in the blanks question to Firestore
# Adds a new "fillTheBlanks" question row (row 17) to the questions sheet,
# mirroring the shape of the existing multipleChoice question rows (9-16),
# and updates the current selection to the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 17: fill-in-the-blanks question -------------------------------
# Row 17 previously held leftover styled-but-empty cells (I17:O17); clear
# that formatting first so the new data row matches the plain style used by
# the other question rows.
$ws.Range("A17:O17").ClearFormats()

$ws.Range("A17").Value = "A1"
$ws.Range("B17").Value = "Reading "
$ws.Range("C17").Value = "Unit1"
$ws.Range("D17").Value = "-"
$ws.Range("E17").Value = "-"
$ws.Range("F17").Value = "What is the meaning of the word"
$ws.Range("G17").Value = "-"
$ws.Range("H17").Value = "Special,اشياء ,bed  ,Home ,Happy "
$ws.Range("I17").Value = "خاص,Things,فراش,منزل ,سعيد"
$ws.Range("J17").Value = "-"
$ws.Range("K17").Value = "fillTheBlanks"
$ws.Range("L17").Value = "-"
$ws.Range("M17").Value = "-"
$ws.Range("N17").Value = "-"
$ws.Range("O17").Value = "-"

# --- Move the active selection to the newly added row ----------------------
$ws.Range("K21").Select()
